$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their literal representation (no numeric auto-coercion)
# by forcing a Text number format before assigning values that look numeric.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.182.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.32"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.44%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.62%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4706"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.64%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.94"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07977"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.006"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.69"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.09%  "

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.984"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.01%  "

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.271"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.37%  "

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.831.88"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.86%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001040"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06590"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.72"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.198.69"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.436"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.08"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.288"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.084.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.04"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.83"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.129"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.476"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.80"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9742"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09496"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.574"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.376"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.344"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02269"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06094"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.420"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.177"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5982"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9994"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1883"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.34"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.275"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5615"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.18"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.980"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06861"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.72"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.75%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.962"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +11.31%  "
